# Updates cryptos list values (price + 1h volume %) to match the latest scrape.
# Numeric-looking "price" strings (column D) are prefixed with a literal leading
# single quote so Excel stores them as text (matching the original inline-string
# cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.985.22"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "'1.639.44"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "'215.03"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'0.5103"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'0.2581"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.06354"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'19.77"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "'0.07759"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'4.278"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "'1.637.34"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "'0.5466"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "'64.35"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "'26.002.02"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'196.47"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'4.427"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "'9.918"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").Value = "'6.083"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "'1.900"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "'0.1225"
$ws.Range("E26").Value = "  +6.80%  "
$ws.Range("D27").Value = "'6.855"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "'15.59"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'1.238"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'0.04856"
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("D31").Value = "'3.277"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'3.210"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'2.374"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").Value = "'0.9149"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").Value = "'2.565"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'1.090.65"
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").Value = "'0.01569"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'2.526"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "'5.579"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("D43").Value = "'0.8049"
$ws.Range("D44").Value = "'99.10"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "0.0₈121"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").Value = "'1.784.06"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'0.4535"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.006"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'55.18"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "'0.05217"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "'7.472"
$ws.Range("E51").Value = "  +0.82%  "
